# Apply the OOXML style changes described by the diff:
#   1. Add a new "Abstract Title" paragraph style (customStyle AbstractTitle),
#      based on Normal, followed by Abstract, centered/bold/colored.
#   2. Change the existing "Abstract" style's spacing-before from 300 to 100
#      (twentieths of a point -> 15pt to 5pt).
#   3. Add a new "Footnote Block Text" paragraph style (styleId
#      FootnoteBlockText), based on Footnote Text, mirroring Block Text's
#      indentation/spacing.

$d = $word.ActiveDocument

# --- 1. New style: Abstract Title ---------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)            # 1 = wdStyleTypeParagraph
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1                   # wdAlignParagraphCenter
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15                # 300 twips

$abstractTitle.Font.Size = 10                                  # sz 20 (half-points)
$abstractTitle.Font.SizeBi = 10                                # szCs 20
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 0x8A5A34                            # OLE BGR -> w:color 345A8A

# --- 2. Abstract: spacing-before 300 -> 100 -----------------------------
$d.Styles("Abstract").ParagraphFormat.SpaceBefore = 5           # 100 twips

# --- 3. New style: Footnote Block Text ----------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)    # 1 = wdStyleTypeParagraph
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceAfter = 5               # 100 twips
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5              # 100 twips
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24              # 480 twips
$footnoteBlockText.ParagraphFormat.RightIndent = 24             # 480 twips

Write-Host "Style edits applied."
